# Generate Report for Handback
# Update "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
# timestamps on the per-language handback status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-25 06:16:44"
$wsZhCn.Range("G5").Value = "2016-02-25 06:17:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-25 06:16:56"
$wsDeDe.Range("G5").Value = "2016-02-25 06:17:50"
